$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (A2 index 0)
$ws.Cells.Item(2, 2).Value = 22.4830047972259
$ws.Cells.Item(2, 3).Value = 4.247072341563142
$ws.Cells.Item(2, 4).Value = 3.649310245373473
$ws.Cells.Item(2, 5).Value = 10.6169397089692
$ws.Cells.Item(2, 6).Value = 57.10392951348858
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 10).Value = 10.48708288890919
$ws.Cells.Item(2, 11).Value = 18.56863169888286
$ws.Cells.Item(2, 12).Value = 11.17439517568099
$ws.Cells.Item(2, 13).Value = 19.43303458844015
$ws.Cells.Item(2, 14).Value = 27.54712734425629
# Row 3 (A3 index 1)
$ws.Cells.Item(3, 2).Value = 22.40308077528776
$ws.Cells.Item(3, 3).Value = 4.108768802822548
$ws.Cells.Item(3, 4).Value = 3.654995034380002
$ws.Cells.Item(3, 5).Value = 10.6323723066171
$ws.Cells.Item(3, 6).Value = 57.08455214805498
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 10).Value = 10.50042270420408
$ws.Cells.Item(3, 11).Value = 18.51658840835018
$ws.Cells.Item(3, 12).Value = 11.19272610962122
$ws.Cells.Item(3, 13).Value = 19.44560374019683
$ws.Cells.Item(3, 14).Value = 27.582061090336
# Row 4 (A4 index 2)
$ws.Cells.Item(4, 2).Value = 22.35904677078684
$ws.Cells.Item(4, 3).Value = 4.022945984096207
$ws.Cells.Item(4, 4).Value = 3.658852132476555
$ws.Cells.Item(4, 5).Value = 10.64254986910862
$ws.Cells.Item(4, 6).Value = 57.0815042544341
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 10).Value = 10.50908855138256
$ws.Cells.Item(4, 11).Value = 18.48870533269229
$ws.Cells.Item(4, 12).Value = 11.20526829870492
$ws.Cells.Item(4, 13).Value = 19.45645224295326
$ws.Cells.Item(4, 14).Value = 27.60535589003678
# Row 5 (A5 index 3)
$ws.Cells.Item(5, 2).Value = 22.34238331292932
$ws.Cells.Item(5, 3).Value = 3.987809609830305
$ws.Cells.Item(5, 4).Value = 3.66051637430521
$ws.Cells.Item(5, 5).Value = 10.64687425184993
$ws.Cells.Item(5, 6).Value = 57.0824899618409
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 10).Value = 10.51273978651496
$ws.Cells.Item(5, 11).Value = 18.47837483562718
$ws.Cells.Item(5, 12).Value = 11.2107033568828
$ws.Cells.Item(5, 13).Value = 19.4616611869301
$ws.Cells.Item(5, 14).Value = 27.61531298627307
# Row 6 (A6 index 4)
$ws.Cells.Item(6, 2).Value = 22.33969407741843
$ws.Cells.Item(6, 3).Value = 3.981967370116277
$ws.Cells.Item(6, 4).Value = 3.660798311230325
$ws.Cells.Item(6, 5).Value = 10.64760301246113
$ws.Cells.Item(6, 6).Value = 57.08278823922335
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 10).Value = 10.51335331926402
$ws.Cells.Item(6, 11).Value = 18.47672201892821
$ws.Cells.Item(6, 12).Value = 11.21162542547521
$ws.Cells.Item(6, 13).Value = 19.46257374545498
$ws.Cells.Item(6, 14).Value = 27.61699440304906
# Row 7 (A7 index 5)
$ws.Cells.Item(7, 2).Value = 22.35881683916851
$ws.Cells.Item(7, 3).Value = 4.022472693884449
$ws.Cells.Item(7, 4).Value = 3.658874202373935
$ws.Cells.Item(7, 5).Value = 10.64260747220663
$ws.Cells.Item(7, 6).Value = 57.08150852564451
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 10).Value = 10.50913730755396
$ws.Cells.Item(7, 11).Value = 18.48856182309478
$ws.Cells.Item(7, 12).Value = 11.20534028535094
$ws.Cells.Item(7, 13).Value = 19.45651930078837
$ws.Cells.Item(7, 14).Value = 27.60548829479879
# Row 8 (A8 index 6)
$ws.Cells.Item(8, 2).Value = 22.45440955662114
$ws.Cells.Item(8, 3).Value = 4.199610172204642
$ws.Cells.Item(8, 4).Value = 3.651194422056675
$ws.Cells.Item(8, 5).Value = 10.62211546177332
$ws.Cells.Item(8, 6).Value = 57.0954130850938
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 10).Value = 10.49158405283791
$ws.Cells.Item(8, 11).Value = 18.54984733232591
$ws.Cells.Item(8, 12).Value = 11.18044882940854
$ws.Cells.Item(8, 13).Value = 19.43671911858447
$ws.Cells.Item(8, 14).Value = 27.55878965041719
# Row 9 (A9 index 7)
$ws.Cells.Item(9, 2).Value = 22.68119387836565
$ws.Cells.Item(9, 3).Value = 4.537202758152089
$ws.Cells.Item(9, 4).Value = 3.639032292060358
$ws.Cells.Item(9, 5).Value = 10.58747990139804
$ws.Cells.Item(9, 6).Value = 57.19276502988237
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 10).Value = 10.46091627062973
$ws.Cells.Item(9, 11).Value = 18.70192784117303
$ws.Cells.Item(9, 12).Value = 11.14182986417551
$ws.Cells.Item(9, 13).Value = 19.42269349141494
$ws.Cells.Item(9, 14).Value = 27.48184492712057
# Row 10 (A10 index 8)
$ws.Cells.Item(10, 2).Value = 22.87080879695531
$ws.Cells.Item(10, 3).Value = 4.776125989529461
$ws.Cells.Item(10, 4).Value = 3.631848833400441
$ws.Cells.Item(10, 5).Value = 10.5653889691543
$ws.Cells.Item(10, 6).Value = 57.30678953642587
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 10).Value = 10.44065103820125
$ws.Cells.Item(10, 11).Value = 18.83249690130078
$ws.Cells.Item(10, 12).Value = 11.11964601149184
$ws.Cells.Item(10, 13).Value = 19.42744507283033
$ws.Cells.Item(10, 14).Value = 27.43421913559763
# Row 11 (A11 index 9)
$ws.Cells.Item(11, 2).Value = 22.96182624180044
$ws.Cells.Item(11, 3).Value = 4.882272588347418
$ws.Cells.Item(11, 4).Value = 3.628958282839019
$ws.Cells.Item(11, 5).Value = 10.55606218919131
$ws.Cells.Item(11, 6).Value = 57.36782288084314
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 10).Value = 10.43191927246381
$ws.Cells.Item(11, 11).Value = 18.89583520365865
$ws.Cells.Item(11, 12).Value = 11.11089250948162
$ws.Cells.Item(11, 13).Value = 19.43285785305987
$ws.Cells.Item(11, 14).Value = 27.4144838215978
# Row 12 (A12 index 10)
$ws.Cells.Item(12, 2).Value = 22.99695483580724
$ws.Cells.Item(12, 3).Value = 4.922061232657058
$ws.Cells.Item(12, 4).Value = 3.627917703416377
$ws.Cells.Item(12, 5).Value = 10.55263382508209
$ws.Cells.Item(12, 6).Value = 57.39224456408935
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 10).Value = 10.42868244918603
$ws.Cells.Item(12, 11).Value = 18.92037236151016
$ws.Cells.Item(12, 12).Value = 11.10776970953698
$ws.Cells.Item(12, 13).Value = 19.43537307265939
$ws.Cells.Item(12, 14).Value = 27.40728788317977
# Row 13 (A13 index 11)
$ws.Cells.Item(13, 2).Value = 22.98936017745362
$ws.Cells.Item(13, 3).Value = 4.913510738938378
$ws.Cells.Item(13, 4).Value = 3.628139412502587
$ws.Cells.Item(13, 5).Value = 10.5533675885953
$ws.Cells.Item(13, 6).Value = 57.38692680669439
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 10).Value = 10.42937646185358
$ws.Cells.Item(13, 11).Value = 18.91506352125078
$ws.Cells.Item(13, 12).Value = 11.10843373018429
$ws.Cells.Item(13, 13).Value = 19.43481070571061
$ws.Cells.Item(13, 14).Value = 27.40882532185493
# Row 14 (A14 index 12)
$ws.Cells.Item(14, 2).Value = 22.96470315144689
$ws.Cells.Item(14, 3).Value = 4.885554397506702
$ws.Cells.Item(14, 4).Value = 3.628871592670559
$ws.Cells.Item(14, 5).Value = 10.55577806395092
$ws.Cells.Item(14, 6).Value = 57.36980586865609
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 10).Value = 10.43165158186742
$ws.Cells.Item(14, 11).Value = 18.89784290731211
$ws.Cells.Item(14, 12).Value = 11.11063175029365
$ws.Cells.Item(14, 13).Value = 19.4330554644127
$ws.Cells.Item(14, 14).Value = 27.41388624865391
# Row 15 (A15 index 13)
$ws.Cells.Item(15, 2).Value = 22.94968557571043
$ws.Cells.Item(15, 3).Value = 4.868376205847166
$ws.Cells.Item(15, 4).Value = 3.62932710075465
$ws.Cells.Item(15, 5).Value = 10.55726801442366
$ws.Cells.Item(15, 6).Value = 57.35948910535774
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 10).Value = 10.43305422696356
$ws.Cells.Item(15, 11).Value = 18.88736626669309
$ws.Cells.Item(15, 12).Value = 11.11200308648654
$ws.Cells.Item(15, 13).Value = 19.43204088869735
$ws.Cells.Item(15, 14).Value = 27.41702233415004
# Row 16 (A16 index 14)
$ws.Cells.Item(16, 2).Value = 22.86495450822506
$ws.Cells.Item(16, 3).Value = 4.769134289412756
$ws.Cells.Item(16, 4).Value = 3.632045308399666
$ws.Cells.Item(16, 5).Value = 10.56601299960221
$ws.Cells.Item(16, 6).Value = 57.30298459094089
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 10).Value = 10.44123145155421
$ws.Cells.Item(16, 11).Value = 18.82843569306784
$ws.Cells.Item(16, 12).Value = 11.12024496045672
$ws.Cells.Item(16, 13).Value = 19.42715659232726
$ws.Cells.Item(16, 14).Value = 27.43554770891713
# Row 17 (A17 index 15)
$ws.Cells.Item(17, 2).Value = 22.81417855363783
$ws.Cells.Item(17, 3).Value = 4.707570923492191
$ws.Cells.Item(17, 4).Value = 3.633809291760197
$ws.Cells.Item(17, 5).Value = 10.57156252143248
$ws.Cells.Item(17, 6).Value = 57.27066243961594
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 10).Value = 10.44637241410685
$ws.Cells.Item(17, 11).Value = 18.79328324019813
$ws.Cells.Item(17, 12).Value = 11.12564347576761
$ws.Cells.Item(17, 13).Value = 19.424991607861
$ws.Cells.Item(17, 14).Value = 27.44740663311784
# Row 18 (A18 index 16)
$ws.Cells.Item(18, 2).Value = 22.78542307548182
$ws.Cells.Item(18, 3).Value = 4.671924332782512
$ws.Cells.Item(18, 4).Value = 3.634859410171592
$ws.Cells.Item(18, 5).Value = 10.57482248864645
$ws.Cells.Item(18, 6).Value = 57.25293460083138
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 10).Value = 10.44937521666559
$ws.Cells.Item(18, 11).Value = 18.77343642405814
$ws.Cells.Item(18, 12).Value = 11.12887452997958
$ws.Cells.Item(18, 13).Value = 19.42405252586683
$ws.Cells.Item(18, 14).Value = 27.45440923060276
# Row 19 (A19 index 17)
$ws.Cells.Item(19, 2).Value = 22.77576482921124
$ws.Cells.Item(19, 3).Value = 4.659815645948461
$ws.Cells.Item(19, 4).Value = 3.6352210703779
$ws.Cells.Item(19, 5).Value = 10.57593795408448
$ws.Cells.Item(19, 6).Value = 57.24708070728821
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 10).Value = 10.45039979982703
$ws.Cells.Item(19, 11).Value = 18.76678094259547
$ws.Cells.Item(19, 12).Value = 11.12999015839165
$ws.Cells.Item(19, 13).Value = 19.42378720936782
$ws.Cells.Item(19, 14).Value = 27.45681139347446
# Row 20 (A20 index 18)
$ws.Cells.Item(20, 2).Value = 22.81953737630215
$ws.Cells.Item(20, 3).Value = 4.714149288824101
$ws.Cells.Item(20, 4).Value = 3.633617838079347
$ws.Cells.Item(20, 5).Value = 10.57096472758646
$ws.Cells.Item(20, 6).Value = 57.27401392398306
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 10).Value = 10.44582040618904
$ws.Cells.Item(20, 11).Value = 18.79698688617799
$ws.Cells.Item(20, 12).Value = 11.12505576025344
$ws.Cells.Item(20, 13).Value = 19.42519040153694
$ws.Cells.Item(20, 14).Value = 27.44612543030611
# Row 21 (A21 index 19)
$ws.Cells.Item(21, 2).Value = 22.97192772746076
$ws.Cells.Item(21, 3).Value = 4.893777199948713
$ws.Cells.Item(21, 4).Value = 3.628655069695177
$ws.Cells.Item(21, 5).Value = 10.5550672436009
$ws.Cells.Item(21, 6).Value = 57.37479922724948
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 10).Value = 10.43098143469494
$ws.Cells.Item(21, 11).Value = 18.90288614922862
$ws.Cells.Item(21, 12).Value = 11.10998093239639
$ws.Cells.Item(21, 13).Value = 19.43355840455376
$ws.Cells.Item(21, 14).Value = 27.41239220425165
# Row 22 (A22 index 20)
$ws.Cells.Item(22, 2).Value = 23.07537268550487
$ws.Cells.Item(22, 3).Value = 5.008787377315198
$ws.Cells.Item(22, 4).Value = 3.625726312293506
$ws.Cells.Item(22, 5).Value = 10.54528032095776
$ws.Cells.Item(22, 6).Value = 57.44829766151437
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 10).Value = 10.42168946679879
$ws.Cells.Item(22, 11).Value = 18.97530944014141
$ws.Cells.Item(22, 12).Value = 11.10124733684462
$ws.Cells.Item(22, 13).Value = 19.44173966430456
$ws.Cells.Item(22, 14).Value = 27.39196241697586
# Row 23 (A23 index 21)
$ws.Cells.Item(23, 2).Value = 23.01981763317807
$ws.Cells.Item(23, 3).Value = 4.947635225136522
$ws.Cells.Item(23, 4).Value = 3.627260728343033
$ws.Cells.Item(23, 5).Value = 10.55044874545306
$ws.Cells.Item(23, 6).Value = 57.40837496243769
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 10).Value = 10.42661170707408
$ws.Cells.Item(23, 11).Value = 18.93636688190577
$ws.Cells.Item(23, 12).Value = 11.10580641516314
$ws.Cells.Item(23, 13).Value = 19.43712572472401
$ws.Cells.Item(23, 14).Value = 27.40271827844914
# Row 24 (A24 index 22)
$ws.Cells.Item(24, 2).Value = 22.81711329183532
$ws.Cells.Item(24, 3).Value = 4.711175994694749
$ws.Cells.Item(24, 4).Value = 3.63370428221567
$ws.Cells.Item(24, 5).Value = 10.57123477376143
$ws.Cells.Item(24, 6).Value = 57.27249605511788
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 10).Value = 10.44606982193322
$ws.Cells.Item(24, 11).Value = 18.79531133662417
$ws.Cells.Item(24, 12).Value = 11.12532106963704
$ws.Cells.Item(24, 13).Value = 19.42509957483503
$ws.Cells.Item(24, 14).Value = 27.44670408653266
# Row 25 (A25 index 23)
$ws.Cells.Item(25, 2).Value = 22.61572791106314
$ws.Cells.Item(25, 3).Value = 4.447265197198652
$ws.Cells.Item(25, 4).Value = 3.642013742490575
$ws.Cells.Item(25, 5).Value = 10.59625847187844
$ws.Cells.Item(25, 6).Value = 57.15894498524472
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 10).Value = 10.46881312225413
$ws.Cells.Item(25, 11).Value = 18.65743003132829
$ws.Cells.Item(25, 12).Value = 11.15118853850928
$ws.Cells.Item(25, 13).Value = 19.42383821501753
$ws.Cells.Item(25, 14).Value = 27.50109561036129
